$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    # Force a numeric-looking string to be stored as text (quote-prefixed
    # entry, like a user typing '602.54 into Excel), then drop the style
    # back to Normal so no visible/number-format change is introduced.
    $cell = $ws.Range($rangeAddr)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "63.656.67"
$ws.Range("E2").Value = "  +0.59%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.313.34"
$ws.Range("E3").Value = "  +4.92%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.19%  "

# Row 5 - BNB
Set-TextValue "D5" "602.54"
$ws.Range("E5").Value = "  +1.82%  "

# Row 6 - Solana
Set-TextValue "D6" "141.96"
$ws.Range("E6").Value = "  +2.15%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.08%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.310.41"
$ws.Range("E8").Value = "  +4.91%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.150"
$ws.Range("E10").Value = "  +2.15%  "

# Row 11 - Toncoin
Set-TextValue "D11" "5.49"
$ws.Range("E11").Value = "  +4.03%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.469"
$ws.Range("E12").Value = "  +1.82%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +0.55%  "

# Row 14 - Avalanche
$ws.Range("E14").Value = "  +0.82%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.858.26"
$ws.Range("E15").Value = "  +4.90%  "

# Row 16 - TRON
$ws.Range("E16").Value = "  +0.13%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.310.01"
$ws.Range("E17").Value = "  +4.95%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "63.737.55"

# Row 19 - Polkadot
Set-TextValue "D19" "6.85"
$ws.Range("E19").Value = "  +2.31%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "479.19"
$ws.Range("E20").Value = "  +0.45%  "

# Row 21 - Chainlink
$ws.Range("E21").Value = "  +0.06%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  +4.21%  "

# Row 23 - Uniswap
Set-TextValue "D23" "8.09"
$ws.Range("E23").Value = "  +4.70%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue "D24" "13.69"
$ws.Range("E24").Value = "  +5.18%  "

# Row 25 - Litecoin
Set-TextValue "D25" "84.49"
$ws.Range("E25").Value = "  -0.18%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.14%  "

# Row 27 - PancakeSwap
$ws.Range("E27").Value = "  +1.25%  "

# Row 28 - NEARProtocol
Set-TextValue "D28" "7.35"
$ws.Range("E28").Value = "  +3.43%  "

# Row 29 - FirstDigitalUSD
$ws.Range("E29").Value = "  -0.29%  "

# Row 30 - RenderToken
Set-TextValue "D30" "8.10"
$ws.Range("E30").Value = "  +0.91%  "

# Row 31 - ImmutableX
$ws.Range("E31").Value = "  +1.21%  "

# Row 32 - EthereumClassic
Set-TextValue "D32" "28.78"
$ws.Range("E32").Value = "  +6.80%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  -0.43%  "

# Row 34 - Stacks
Set-TextValue "D34" "2.53"
$ws.Range("E34").Value = "  -0.51%  "

# Row 35 - Mantle
$ws.Range("E35").Value = "  +2.54%  "

# Row 36 - Filecoin
Set-TextValue "D36" "6.01"
$ws.Range("E36").Value = "  +3.18%  "

# Row 37 - OKB
Set-TextValue "D37" "52.54"
$ws.Range("E37").Value = "  -0.01%  "

# Row 38 - PEPE
$ws.Range("D38").Value = "0.0₃0738"
$ws.Range("E38").Value = "  +4.78%  "

# Row 39 - VeChain
Set-TextValue "D39" "0.0400"
$ws.Range("E39").Value = "  +2.32%  "

# Row 40 - Bittensor
Set-TextValue "D40" "432.62"
$ws.Range("E40").Value = "  +2.34%  "

# Row 41 - Maker
$ws.Range("D41").Value = "3.057.34"
$ws.Range("E41").Value = "  +3.84%  "

# Row 42 - now dogwifhat (was Cosmos)
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D42" "2.74"
$ws.Range("E42").Value = "  -0.75%  "

# Row 43 - now Cosmos (was dogwifhat)
$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D43" "8.32"
$ws.Range("E43").Value = "  +0.38%  "

# Row 44 - Kaspa
$ws.Range("E44").Value = "  +2.87%  "

# Row 45 - TheGraph
$ws.Range("E45").Value = "  -0.44%  "

# Row 46 - Fetch.AI
$ws.Range("E46").Value = "  +2.67%  "

# Row 47 - InjectiveProtocol
Set-TextValue "D47" "26.27"
$ws.Range("E47").Value = "  +2.77%  "

# Row 48 - Arweave
Set-TextValue "D48" "36.06"
$ws.Range("E48").Value = "  +12.77%  "

# Row 49 - USDe
$ws.Range("E49").Value = "  -0.02%  "

# Row 50 - Monero
Set-TextValue "D50" "125.21"
$ws.Range("E50").Value = "  +3.46%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  +0.71%  "
